$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1042.4286
$ws.Range("I38").Value = 324.25
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 972.75
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -600.75
$ws.Range("N38").Value = -6744

$ws.Range("H58").Value = 1527.5714
$ws.Range("I58").Value = 448.25
$ws.Range("J58").Value = 2966.6667
$ws.Range("K58").Value = 1344.75
$ws.Range("L58").Value = 8900.000100000001
$ws.Range("M58").Value = -1194.75
$ws.Range("N58").Value = -9200.000100000001

$ws.Range("H82").Value = 1786.3334
$ws.Range("I82").Value = 679.5
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 2038.5
$ws.Range("L82").Value = 12000
$ws.Range("M82").Value = -1632.5
$ws.Range("N82").Value = -12812

$ws.Range("H85").Value = 1786.3334
$ws.Range("I85").Value = 679.5
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 2038.5
$ws.Range("L85").Value = 12000
$ws.Range("M85").Value = -634.5
$ws.Range("N85").Value = -14808

$ws.Range("H87").Value = 18833.658
$ws.Range("J87").Value = 18833.658
$ws.Range("L87").Value = 18833.658
$ws.Range("N87").Value = -21329.658

$ws.Range("H90").Value = 18833.658
$ws.Range("J90").Value = 18833.658
$ws.Range("L90").Value = 56500.974
$ws.Range("N90").Value = -68980.974

$ws.Range("H129").Value = 869.51666
$ws.Range("J129").Value = 870.9828
$ws.Range("L129").Value = 2612.9484
$ws.Range("N129").Value = -12612.9484

$ws.Range("H137").Value = 1932.6666
$ws.Range("I137").Value = 1405.4
$ws.Range("J137").Value = 2811.4443
$ws.Range("K137").Value = 4216.200000000001
$ws.Range("L137").Value = 8434.332900000001
$ws.Range("M137").Value = -1666.200000000001
$ws.Range("N137").Value = -13534.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6548.478
$ws.Range("I32").Value = 4429.2
$ws.Range("K32").Value = 4429.2
$ws.Range("M32").Value = -4142.2

$ws.Range("H61").Value = 4299.5654
$ws.Range("I61").Value = 3325.6667
$ws.Range("J61").Value = 7805.6
$ws.Range("K61").Value = 3325.6667
$ws.Range("L61").Value = 7805.6
$ws.Range("M61").Value = -3113.6667
$ws.Range("N61").Value = -8229.6

$ws.Range("H97").Value = 1220.1316
$ws.Range("I97").Value = 1007.9375
$ws.Range("J97").Value = 2351.8333
$ws.Range("K97").Value = 1007.9375
$ws.Range("L97").Value = 2351.8333
$ws.Range("M97").Value = -511.9375
$ws.Range("N97").Value = -3343.8333

$ws.Range("H122").Value = 2474.9473
$ws.Range("I122").Value = 1729.2858
$ws.Range("J122").Value = 4562.8
$ws.Range("K122").Value = 5187.857400000001
$ws.Range("L122").Value = 13688.4
$ws.Range("M122").Value = -2737.857400000001
$ws.Range("N122").Value = -18588.4

$ws.Range("H132").Value = 3146
$ws.Range("I132").Value = 1351.3334
$ws.Range("J132").Value = 5593.273
$ws.Range("K132").Value = 4054.0002
$ws.Range("L132").Value = 16779.819
$ws.Range("M132").Value = -1524.0002
$ws.Range("N132").Value = -21839.819

$ws.Range("H136").Value = 4299.5654
$ws.Range("I136").Value = 3325.6667
$ws.Range("J136").Value = 7805.6
$ws.Range("K136").Value = 9977.000100000001
$ws.Range("L136").Value = 23416.8
$ws.Range("M136").Value = -7427.000100000001
$ws.Range("N136").Value = -28516.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1278.72
$ws.Range("I99").Value = 879.3125
$ws.Range("J99").Value = 1988.7778
$ws.Range("K99").Value = 879.3125
$ws.Range("L99").Value = 1988.7778
$ws.Range("M99").Value = 618.6875
$ws.Range("N99").Value = -4984.7778

$ws.Range("H134").Value = 2631.3333
$ws.Range("I134").Value = 1314.6
$ws.Range("J134").Value = 3948.0667
$ws.Range("K134").Value = 3943.8
$ws.Range("L134").Value = 11844.2001
$ws.Range("M134").Value = -1408.8
$ws.Range("N134").Value = -16914.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2017.9
$ws.Range("I16").Value = 2019.8889
$ws.Range("K16").Value = 2019.8889
$ws.Range("M16").Value = -1732.8889

$ws.Range("H58").Value = 2427.5334
$ws.Range("I58").Value = 1950.5
$ws.Range("J58").Value = 2745.5557
$ws.Range("K58").Value = 1950.5
$ws.Range("L58").Value = 2745.5557
$ws.Range("M58").Value = -1747.5
$ws.Range("N58").Value = -3151.5557

$ws.Range("H113").Value = 2017.9
$ws.Range("I113").Value = 2019.8889
$ws.Range("K113").Value = 2019.8889
$ws.Range("M113").Value = 150.1111000000001

$ws.Range("H132").Value = 2469.5652
$ws.Range("I132").Value = 1421.2727
$ws.Range("J132").Value = 3430.5
$ws.Range("K132").Value = 4263.8181
$ws.Range("L132").Value = 10291.5
$ws.Range("M132").Value = -1733.8181
$ws.Range("N132").Value = -15351.5

$ws.Range("H134").Value = 3003
$ws.Range("I134").Value = 2727.7222
$ws.Range("K134").Value = 8183.1666
$ws.Range("M134").Value = -5648.1666

$ws.Range("H136").Value = 2427.5334
$ws.Range("I136").Value = 1950.5
$ws.Range("J136").Value = 2745.5557
$ws.Range("K136").Value = 5851.5
$ws.Range("L136").Value = 8236.667099999999
$ws.Range("M136").Value = -3301.5
$ws.Range("N136").Value = -13336.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3710
$ws.Range("I56").Value = 3710
$ws.Range("K56").Value = 3710
$ws.Range("M56").Value = -3180

$ws.Range("H107").Value = 654.2727
$ws.Range("I107").Value = 708.8333
$ws.Range("J107").Value = 588.8
$ws.Range("K107").Value = 2126.4999
$ws.Range("L107").Value = 1766.4
$ws.Range("M107").Value = -206.4998999999998
$ws.Range("N107").Value = -5606.4

$ws.Range("H131").Value = 2340.6172
$ws.Range("J131").Value = 2495.946
$ws.Range("L131").Value = 7487.838
$ws.Range("N131").Value = -17567.838

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4175.1904
$ws.Range("I132").Value = 3071.1333
$ws.Range("J132").Value = 6935.3335
$ws.Range("K132").Value = 9213.3999
$ws.Range("L132").Value = 20806.0005
$ws.Range("M132").Value = -6683.3999
$ws.Range("N132").Value = -25866.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1895.7391
$ws.Range("I22").Value = 878.7143
$ws.Range("J22").Value = 2340.6875
$ws.Range("K22").Value = 878.7143
$ws.Range("L22").Value = 2340.6875
$ws.Range("M22").Value = -583.7143
$ws.Range("N22").Value = -2930.6875

$ws.Range("H27").Value = 1895.7391
$ws.Range("I27").Value = 878.7143
$ws.Range("J27").Value = 2340.6875
$ws.Range("K27").Value = 878.7143
$ws.Range("L27").Value = 2340.6875
$ws.Range("M27").Value = -771.7143
$ws.Range("N27").Value = -2554.6875

$ws.Range("H82").Value = 3725
$ws.Range("I82").Value = 2950
$ws.Range("J82").Value = 4500
$ws.Range("K82").Value = 2950
$ws.Range("L82").Value = 4500
$ws.Range("M82").Value = -2589
$ws.Range("N82").Value = -5222

$ws.Range("H85").Value = 3725
$ws.Range("I85").Value = 2950
$ws.Range("J85").Value = 4500
$ws.Range("K85").Value = 2950
$ws.Range("L85").Value = 4500
$ws.Range("M85").Value = -1702
$ws.Range("N85").Value = -6996

$ws.Range("H132").Value = 14818.565
$ws.Range("I132").Value = 17551.188
$ws.Range("J132").Value = 8572.571
$ws.Range("K132").Value = 52653.564
$ws.Range("L132").Value = 25717.713
$ws.Range("M132").Value = -50123.564
$ws.Range("N132").Value = -30777.713

$ws.Range("H136").Value = 19613568
$ws.Range("I136").Value = 3782
$ws.Range("J136").Value = 47627548
$ws.Range("K136").Value = 11346
$ws.Range("L136").Value = 142882644
$ws.Range("M136").Value = -8796
$ws.Range("N136").Value = -142887744

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3205
$ws.Range("I132").Value = 2583.5
$ws.Range("K132").Value = 7750.5
$ws.Range("M132").Value = -5220.5

$ws.Range("H136").Value = 2977.2693
$ws.Range("I136").Value = 1335.8182
$ws.Range("J136").Value = 4181
$ws.Range("K136").Value = 4007.4546
$ws.Range("L136").Value = 12543
$ws.Range("M136").Value = -1457.4546
$ws.Range("N136").Value = -17643
